# Applies the "Ajout du diagramme de classe" commit to Journal_de_travail.xlsx
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuil1")
$ws2 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Fill in the new journal entries in Tableau1 (rows 13-22 on Feuil1).
# Numeric / date cells first (they never touch the shared-string table), the
# text cells are written afterwards in the exact order the original author
# used so the shared-string table comes out in the same order as the diff.
# ---------------------------------------------------------------------------

# Row 13
$ws1.Range("B13").Value = 45412
$ws1.Range("C13").Value = 0.66666666666666663
$ws1.Range("D13").Value = 0.69444444444444453
$ws1.Range("F13").Value = "Analyse"

# Row 14
$ws1.Range("B14").Value = 45412
$ws1.Range("C14").Value = 0.69444444444444453
$ws1.Range("D14").Value = 0.70486111111111116
$ws1.Range("F14").Value = "Autres"

# Row 15
$ws1.Range("B15").Value = 45414
$ws1.Range("C15").Value = 0.36805555555555558
$ws1.Range("D15").Value = 0.39930555555555558
$ws1.Range("F15").Value = "Analyse"

# Row 16
$ws1.Range("B16").Value = 45414
$ws1.Range("C16").Value = 0.40972222222222227
$ws1.Range("D16").Value = 0.4513888888888889
$ws1.Range("F16").Value = "Analyse"

# Row 17
$ws1.Range("B17").Value = 45414
$ws1.Range("C17").Value = 0.4513888888888889
$ws1.Range("D17").Value = 0.51041666666666663
$ws1.Range("F17").Value = "Analyse"

# Row 18
$ws1.Range("B18").Value = 45414
$ws1.Range("C18").Value = 0.5625
$ws1.Range("D18").Value = 0.59375
$ws1.Range("F18").Value = "Analyse"

# Row 19
$ws1.Range("B19").Value = 45414
$ws1.Range("C19").Value = 0.59375
$ws1.Range("D19").Value = 0.60416666666666663
$ws1.Range("F19").Value = "Analyse"

# Row 20
$ws1.Range("B20").Value = 45414
$ws1.Range("C20").Value = 0.61458333333333337
$ws1.Range("D20").Value = 0.62847222222222221
$ws1.Range("F20").Value = "Analyse"

# Row 21
$ws1.Range("B21").Value = 45414
$ws1.Range("C21").Value = 0.63888888888888895
$ws1.Range("D21").Value = 0.64930555555555558
$ws1.Range("F21").Value = "Documentation"

# Row 22
$ws1.Range("B22").Value = 45414
$ws1.Range("C22").Value = 0.64930555555555558
$ws1.Range("D22").Value = 0.67013888888888884
$ws1.Range("F22").Value = "Analyse"

# Text cells (Description / Remarque), in the original authoring order
$ws1.Range("G14").Value = "Écriture de l'email pour l'envoi de la planification initiale"
$ws1.Range("G17").Value = "Création du diagramme de classes"
$ws1.Range("H14").Value = "Je me sens confient pour la suite du projet "
$ws1.Range("G15").Value = "Création des milestones, des tags et du projet KanBan sur GitHub"
$ws1.Range("G16").Value = "Création des tâches sur GitHub"
$ws1.Range("G18").Value = "Réalisation du diagramme de classes"
$ws1.Range("G19").Value = "Validation du diagramme de classes"
$ws1.Range("I19").Value = "Validation faite avec Monsieur Viret"
$ws1.Range("G20").Value = "Rectification du diagramme de classes"
$ws1.Range("G21").Value = "Écriture de l'analyse du diagramme de classes"
$ws1.Range("G13").Value = "Écriture de la planification initiale"
$ws1.Range("G22").Value = "Création du diagramme de flux la gestion du mouvement des pièces"

# Hyperlinked notes (row 18 then row 17, matching the original order), each
# cell keeps its two-line "title + URL" text while the hyperlink's display
# text is just the bare URL.
$ws1.Hyperlinks.Add($ws1.Range("I18"), "https://cpnv-es-ngy.gitbook.io/uml-backlog/class-diagram/standards/les-relations", "", "", "https://cpnv-es-ngy.gitbook.io/uml-backlog/class-diagram/standards/les-relations")
$ws1.Range("I18").Value = "Guide pour les relations:" + [char]10 + "https://cpnv-es-ngy.gitbook.io/uml-backlog/class-diagram/standards/les-relations"

$ws1.Hyperlinks.Add($ws1.Range("I17"), "https://astah.net/support/astah-pro/user-guide/class-diagrams/", "", "", "https://astah.net/support/astah-pro/user-guide/class-diagrams/")
$ws1.Range("I17").Value = "Comment utiliser astah:" + [char]10 + "https://astah.net/support/astah-pro/user-guide/class-diagrams/"

# Rows 17/18/22 grow taller to host the wrapped two-/three-line text
$ws1.Rows.Item(17).RowHeight = 30
$ws1.Rows.Item(18).RowHeight = 45
$ws1.Rows.Item(22).RowHeight = 30

# Widen column I a bit to best-fit the new, longer hyperlink notes
$ws1.Columns.Item(9).ColumnWidth = 58.95

# ---------------------------------------------------------------------------
# View state: the workbook was left with Sheet1 active (2nd tab), scrolled
# down a bit on Feuil1, with new selections on both sheets.
# ---------------------------------------------------------------------------
$ws1.Range("G23").Select()
$ws2.Range("K20").Select()
$ws2.Activate()
